# The commit swaps the contents of ppt/theme/theme1.xml (the theme used
# by the presentation's single slide master, i.e. the theme every slide
# actually renders with) and ppt/theme/theme2.xml (a theme used only by
# the notes master, which PowerPoint's object model does not expose for
# editing) - in effect, the author applied the built-in "Office Theme"
# design (Design tab) in place of the previous "Integral" design.
#
# The font scheme (majorFont/minorFont) and the format scheme
# (fillStyleLst/lnStyleLst/effectStyleLst/bgFillStyleLst) are byte-for-byte
# identical between the old "Integral" theme and the new "Office Theme",
# so the only observable difference applying the new design makes to the
# live theme part is its 12-slot colour scheme. We reproduce that by
# writing the "Office" theme colours into the active theme's colour
# scheme via ThemeColorScheme, which PowerPoint re-serialises back into
# ppt/theme/theme1.xml's <a:clrScheme> on save.
#
# PowerPoint colour properties use the OLE "BGR" integer encoding
# (0x00BBGGRR), i.e. R | (G << 8) | (B << 16).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
